# Append a new logbook entry (row 53) continuing the "10:50AM" start time,
# plus a new row 54 that begins the next entry at "2:00PM" (reusing the
# existing shared string), matching the "newProduct logic finished" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C53").Value = "12:40PM"
$ws.Range("D53").Value = 10
$ws.Range("E53").Value = 100
$ws.Range("F53").Value = "Finalizing searchEditClients page/creating and adding newProduct page"
$ws.Range("G53").Value = "newProduct page now just needs the code to add prods to db."

$ws.Range("B54").Value = "2:00PM"

# Update the active selection to follow the newly entered data, as Excel
# would after the user finished typing in B54.
$ws.Range("B54").Select()
